# Auto-generated edit script: updates numeric cells in Leve profit sheets
# to reflect refreshed market-board prices (per scheduled runner update).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3624.7778
$ws.Range("I40").Value = 3191
$ws.Range("J40").Value = 3923
$ws.Range("K40").Value = 3191
$ws.Range("L40").Value = 3923
$ws.Range("M40").Value = -3016
$ws.Range("N40").Value = -4273

# Row 43
$ws.Range("H43").Value = 3423.2
$ws.Range("I43").Value = 1135.1
$ws.Range("J43").Value = 7999.4
$ws.Range("K43").Value = 1135.1
$ws.Range("L43").Value = 7999.4
$ws.Range("M43").Value = -1066.1
$ws.Range("N43").Value = -8137.4


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2166.261
$ws.Range("I45").Value = 1846.5454
$ws.Range("J45").Value = 2459.3333
$ws.Range("K45").Value = 1846.5454
$ws.Range("L45").Value = 2459.3333
$ws.Range("M45").Value = -1469.5454
$ws.Range("N45").Value = -3213.3333

# Row 61
$ws.Range("H61").Value = 55674180
$ws.Range("I61").Value = 250001760
$ws.Range("K61").Value = 250001760
$ws.Range("M61").Value = -250001548

# Row 110
$ws.Range("H110").Value = 4022
$ws.Range("I110").Value = 3950.875
$ws.Range("J110").Value = 4306.5
$ws.Range("K110").Value = 3950.875
$ws.Range("L110").Value = 4306.5
$ws.Range("M110").Value = -1905.875
$ws.Range("N110").Value = -8396.5

# Row 132
$ws.Range("H132").Value = 4792.027
$ws.Range("I132").Value = 2371.8518
$ws.Range("K132").Value = 7115.555399999999
$ws.Range("M132").Value = -4585.555399999999

# Row 136
$ws.Range("H136").Value = 55674180
$ws.Range("I136").Value = 250001760
$ws.Range("K136").Value = 750005280
$ws.Range("M136").Value = -750002730


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 32997.145
$ws.Range("I82").Value = 10330
$ws.Range("J82").Value = 49997.5
$ws.Range("K82").Value = 10330
$ws.Range("L82").Value = 49997.5
$ws.Range("M82").Value = -9947
$ws.Range("N82").Value = -50763.5

# Row 85
$ws.Range("H85").Value = 32997.145
$ws.Range("I85").Value = 10330
$ws.Range("J85").Value = 49997.5
$ws.Range("K85").Value = 10330
$ws.Range("L85").Value = 49997.5
$ws.Range("M85").Value = -9004
$ws.Range("N85").Value = -52649.5

# Row 134
$ws.Range("H134").Value = 36534.594
$ws.Range("I134").Value = 5261.8623
$ws.Range("K134").Value = 15785.5869
$ws.Range("M134").Value = -13250.5869


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4174.4287
$ws.Range("I99").Value = 3880.5
$ws.Range("K99").Value = 3880.5
$ws.Range("M99").Value = -2382.5

# Row 107
$ws.Range("H107").Value = 1433.2667
$ws.Range("I107").Value = 1146.875
$ws.Range("K107").Value = 1146.875
$ws.Range("M107").Value = 773.125

# Row 117
$ws.Range("H117").Value = 39944.5
$ws.Range("J117").Value = 39944.5
$ws.Range("L117").Value = 39944.5
$ws.Range("N117").Value = -49122.5

# Row 126
$ws.Range("H126").Value = 4174.4287
$ws.Range("I126").Value = 3880.5
$ws.Range("K126").Value = 11641.5
$ws.Range("M126").Value = -9171.5


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 396.27274
$ws.Range("I11").Value = 295.44446
$ws.Range("K11").Value = 886.33338
$ws.Range("M11").Value = -746.33338

# Row 18
$ws.Range("H18").Value = 2500
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 7500
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -7838

# Row 47
$ws.Range("H47").Value = 359.66666
$ws.Range("I47").Value = 359.66666
$ws.Range("K47").Value = 1078.99998
$ws.Range("M47").Value = -647.9999800000001

# Row 52
$ws.Range("H52").Value = 12000
$ws.Range("J52").Value = 12000
$ws.Range("L52").Value = 36000
$ws.Range("N52").Value = -36532

# Row 118
$ws.Range("H118").Value = 6633.143
$ws.Range("I118").Value = 4800
$ws.Range("J118").Value = 8008
$ws.Range("K118").Value = 14400
$ws.Range("L118").Value = 24024
$ws.Range("M118").Value = -13157
$ws.Range("N118").Value = -26510

# Row 131
$ws.Range("H131").Value = 3681.8462
$ws.Range("I131").Value = 1415.091
$ws.Range("K131").Value = 4245.272999999999
$ws.Range("M131").Value = 794.7270000000008


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3882.4
$ws.Range("I113").Value = 2937
$ws.Range("K113").Value = 2937
$ws.Range("M113").Value = -767

# Row 118
$ws.Range("H118").Value = 42550
$ws.Range("J118").Value = 42550
$ws.Range("L118").Value = 42550
$ws.Range("N118").Value = -45864

# Row 130
$ws.Range("H130").Value = 85999.664
$ws.Range("J130").Value = 85999.664
$ws.Range("L130").Value = 85999.664
$ws.Range("N130").Value = -96039.664

# Row 132
$ws.Range("H132").Value = 28574514
$ws.Range("I132").Value = 32261002
$ws.Range("K132").Value = 96783006
$ws.Range("M132").Value = -96780476


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4830

# Row 22
$ws.Range("H22").Value = 1584.5
$ws.Range("I22").Value = 1809.5
$ws.Range("K22").Value = 1809.5
$ws.Range("M22").Value = -1514.5

# Row 27
$ws.Range("H27").Value = 1584.5
$ws.Range("I27").Value = 1809.5
$ws.Range("K27").Value = 1809.5
$ws.Range("M27").Value = -1702.5

# Row 46
$ws.Range("H46").Value = 3134.0715
$ws.Range("I46").Value = 3070.2666
$ws.Range("K46").Value = 3070.2666
$ws.Range("M46").Value = -2882.2666

# Row 68
$ws.Range("H68").Value = 4925.5
$ws.Range("I68").Value = 1900.6666
$ws.Range("J68").Value = 14000
$ws.Range("K68").Value = 1900.6666
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = -1151.6666
$ws.Range("N68").Value = -15498

# Row 71
$ws.Range("H71").Value = 4925.5
$ws.Range("I71").Value = 1900.6666
$ws.Range("J71").Value = 14000
$ws.Range("K71").Value = 9503.333000000001
$ws.Range("L71").Value = 70000
$ws.Range("M71").Value = -5759.333000000001
$ws.Range("N71").Value = -77488


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 41296
$ws.Range("I81").Value = 1660
$ws.Range("J81").Value = 100750
$ws.Range("K81").Value = 3320
$ws.Range("L81").Value = 201500
$ws.Range("M81").Value = -2259
$ws.Range("N81").Value = -203622

# Row 84
$ws.Range("H84").Value = 41296
$ws.Range("I84").Value = 1660
$ws.Range("J84").Value = 100750
$ws.Range("K84").Value = 16600
$ws.Range("L84").Value = 1007500
$ws.Range("M84").Value = -11296
$ws.Range("N84").Value = -1018108

# Row 96
$ws.Range("H96").Value = 1979.6666
$ws.Range("I96").Value = 1989
$ws.Range("J96").Value = 1975
$ws.Range("K96").Value = 1989
$ws.Range("L96").Value = 1975
$ws.Range("M96").Value = -616
$ws.Range("N96").Value = -4721

# Row 100
$ws.Range("H100").Value = 2007.151
$ws.Range("I100").Value = 1962.5116
$ws.Range("K100").Value = 3925.0232
$ws.Range("M100").Value = -3384.0232

# Row 126
$ws.Range("H126").Value = 4730.3335
$ws.Range("I126").Value = 4901.778
$ws.Range("K126").Value = 14705.334
$ws.Range("M126").Value = -12235.334

# Row 136
$ws.Range("H136").Value = 1637.65
$ws.Range("I136").Value = 1585.0667
$ws.Range("J136").Value = 1795.4
$ws.Range("K136").Value = 4755.2001
$ws.Range("L136").Value = 5386.200000000001
$ws.Range("M136").Value = -2205.2001
$ws.Range("N136").Value = -10486.2

